$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D9","D10","D11","D12","D14","D16","D19","D20","D21","D22","D23","D24","D26","D27","D28","D30","D31","D32","D33","D34","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '44.247.43'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '2.290.90'
$ws.Range("E3").Value = '  +3.93%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '274.59'
$ws.Range("E5").Value = '  +5.88%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '96.46'
$ws.Range("E6").Value = '  +12.07%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +1.64%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.638'
$ws.Range("E9").Value = '  +7.28%  '
$ws.Range("D10").Value = '47.64'
$ws.Range("E10").Value = '  +5.58%  '
$ws.Range("D11").Value = '0.0941'
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("D12").Value = '8.16'
$ws.Range("E12").Value = '  +10.98%  '
$ws.Range("D14").Value = '15.72'
$ws.Range("E14").Value = '  +9.14%  '
$ws.Range("D15").Value = '2.623.71'
$ws.Range("E15").Value = '  +3.49%  '
$ws.Range("D16").Value = '0.841'
$ws.Range("E16").Value = '  +8.29%  '
$ws.Range("D17").Value = '2.291.27'
$ws.Range("E17").Value = '  +4.09%  '
$ws.Range("D18").Value = '44.274.40'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").Value = '0.0000106'
$ws.Range("E19").Value = '  +3.04%  '
$ws.Range("D20").Value = '6.24'
$ws.Range("E20").Value = '  +5.59%  '
$ws.Range("D21").Value = '71.17'
$ws.Range("E21").Value = '  +1.92%  '
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").Value = '10.30'
$ws.Range("E22").Value = '  +15.69%  '
$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").Value = '2.32'
$ws.Range("E23").Value = '  -0.56%  '
$ws.Range("D24").Value = '235.90'
$ws.Range("E24").Value = '  +2.13%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '11.53'
$ws.Range("E26").Value = '  +8.34%  '
$ws.Range("D27").Value = '2.52'
$ws.Range("E27").Value = '  +12.40%  '
$ws.Range("D28").Value = '39.69'
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("E29").Value = '  -5.54%  '
$ws.Range("D30").Value = '2.26'
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").Value = '173.71'
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '22.04'
$ws.Range("E32").Value = '  +8.16%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0921'
$ws.Range("E33").Value = '  +6.21%  '
$ws.Range("D34").Value = '5.69'
$ws.Range("E34").Value = '  +6.47%  '
$ws.Range("E35").Value = '  +1.88%  '
$ws.Range("E36").Value = '  +4.10%  '
$ws.Range("D37").Value = '0.0357'
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").Value = '4.45'
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").Value = '3.55'
$ws.Range("E39").Value = '  +24.64%  '
$ws.Range("D40").Value = '0.250'
$ws.Range("E40").Value = '  +26.12%  '
$ws.Range("D41").Value = '2.23'
$ws.Range("E41").Value = '  +6.78%  '
$ws.Range("D42").Value = '12.59'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").Value = '5.54'
$ws.Range("E43").Value = '  +1.23%  '
$ws.Range("D44").Value = '62.75'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("E45").Value = '  +5.31%  '
$ws.Range("D46").Value = '8.62'
$ws.Range("E46").Value = '  +3.35%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '1.19'
$ws.Range("E47").Value = '  +7.55%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '100.82'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("D51").Value = '2.508.03'
$ws.Range("E51").Value = '  +3.54%  '
